$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 456  # H2
$ws.Cells.Item(2, 9).Value = 456  # I2
$ws.Cells.Item(2, 11).Value = 456  # K2
$ws.Cells.Item(2, 13).Value = -343  # M2

$ws.Cells.Item(70, 8).Value = 2931.6667  # H70
$ws.Cells.Item(70, 10).Value = 3178.2  # J70
$ws.Cells.Item(70, 12).Value = 9534.599999999999  # L70
$ws.Cells.Item(70, 14).Value = -10074.6  # N70

$ws.Cells.Item(73, 8).Value = 2931.6667  # H73
$ws.Cells.Item(73, 10).Value = 3178.2  # J73
$ws.Cells.Item(73, 12).Value = 9534.599999999999  # L73
$ws.Cells.Item(73, 14).Value = -11406.6  # N73

$ws.Cells.Item(74, 8).Value = 4611.4443  # H74
$ws.Cells.Item(74, 9).Value = 4562.875  # I74
$ws.Cells.Item(74, 11).Value = 4562.875  # K74
$ws.Cells.Item(74, 13).Value = -3626.875  # M74

$ws.Cells.Item(77, 8).Value = 4611.4443  # H77
$ws.Cells.Item(77, 9).Value = 4562.875  # I77
$ws.Cells.Item(77, 11).Value = 22814.375  # K77
$ws.Cells.Item(77, 13).Value = -18134.375  # M77

$ws.Cells.Item(132, 8).Value = 4037.0557  # H132
$ws.Cells.Item(132, 9).Value = 2024.2727  # I132
$ws.Cells.Item(132, 11).Value = 6072.8181  # K132
$ws.Cells.Item(132, 13).Value = -3542.8181  # M132


# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 12487.889  # H122
$ws.Cells.Item(122, 9).Value = 11770.286  # I122
$ws.Cells.Item(122, 10).Value = 14999.5  # J122
$ws.Cells.Item(122, 11).Value = 35310.858  # K122
$ws.Cells.Item(122, 12).Value = 44998.5  # L122
$ws.Cells.Item(122, 13).Value = -32860.858  # M122
$ws.Cells.Item(122, 14).Value = -49898.5  # N122


# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 7500  # H26
$ws.Cells.Item(26, 9).Value = 7500  # I26
$ws.Cells.Item(26, 11).Value = 7500  # K26
$ws.Cells.Item(26, 13).Value = -7208  # M26

$ws.Cells.Item(96, 8).Value = 5646.6665  # H96
$ws.Cells.Item(96, 9).Value = 5646.6665  # I96
$ws.Cells.Item(96, 11).Value = 5646.6665  # K96
$ws.Cells.Item(96, 13).Value = -2900.6665  # M96


# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 98.5  # H7
$ws.Cells.Item(7, 9).Value = 189.5  # I7
$ws.Cells.Item(7, 10).Value = 7.5  # J7
$ws.Cells.Item(7, 11).Value = 189.5  # K7
$ws.Cells.Item(7, 12).Value = 7.5  # L7
$ws.Cells.Item(7, 13).Value = -76.5  # M7
$ws.Cells.Item(7, 14).Value = -233.5  # N7

$ws.Cells.Item(16, 8).Value = 840.25  # H16
$ws.Cells.Item(16, 10).Value = 0  # J16
$ws.Cells.Item(16, 12).Value = 0  # L16
$ws.Cells.Item(16, 14).ClearContents()  # N16

$ws.Cells.Item(68, 8).Value = 50000  # H68
$ws.Cells.Item(68, 10).Value = 50000  # J68
$ws.Cells.Item(68, 12).Value = 50000  # L68
$ws.Cells.Item(68, 14).Value = -51498  # N68

$ws.Cells.Item(71, 8).Value = 50000  # H71
$ws.Cells.Item(71, 10).Value = 50000  # J71
$ws.Cells.Item(71, 12).Value = 150000  # L71
$ws.Cells.Item(71, 14).Value = -157488  # N71

$ws.Cells.Item(74, 8).Value = 30000  # H74
$ws.Cells.Item(74, 9).Value = 30000  # I74
$ws.Cells.Item(74, 11).Value = 30000  # K74
$ws.Cells.Item(74, 13).Value = -29126  # M74

$ws.Cells.Item(77, 8).Value = 30000  # H77
$ws.Cells.Item(77, 9).Value = 30000  # I77
$ws.Cells.Item(77, 11).Value = 90000  # K77
$ws.Cells.Item(77, 13).Value = -85632  # M77

$ws.Cells.Item(93, 8).Value = 3800.6667  # H93
$ws.Cells.Item(93, 9).Value = 3800.6667  # I93
$ws.Cells.Item(93, 11).Value = 3800.6667  # K93
$ws.Cells.Item(93, 13).Value = -1928.6667  # M93

$ws.Cells.Item(113, 8).Value = 840.25  # H113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 14).ClearContents()  # N113

$ws.Cells.Item(132, 8).Value = 2955.4546  # H132
$ws.Cells.Item(132, 9).Value = 2667.7778  # I132
$ws.Cells.Item(132, 11).Value = 8003.3334  # K132
$ws.Cells.Item(132, 13).Value = -5473.3334  # M132


# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 37.916668  # H12
$ws.Cells.Item(12, 10).Value = 30.125  # J12
$ws.Cells.Item(12, 12).Value = 90.375  # L12
$ws.Cells.Item(12, 14).Value = -436.375  # N12

$ws.Cells.Item(132, 8).Value = 1823.75  # H132
$ws.Cells.Item(132, 9).Value = 1450.25  # I132
$ws.Cells.Item(132, 10).Value = 2197.25  # J132
$ws.Cells.Item(132, 11).Value = 13052.25  # K132
$ws.Cells.Item(132, 12).Value = 19775.25  # L132
$ws.Cells.Item(132, 13).Value = -10522.25  # M132
$ws.Cells.Item(132, 14).Value = -24835.25  # N132

$ws.Cells.Item(133, 8).Value = 3076.6667  # H133
$ws.Cells.Item(133, 9).Value = 3076.6667  # I133
$ws.Cells.Item(133, 11).Value = 9230.000100000001  # K133
$ws.Cells.Item(133, 13).Value = -4170.000100000001  # M133

$ws.Cells.Item(134, 8).Value = 997.5  # H134
$ws.Cells.Item(134, 9).Value = 997.5  # I134
$ws.Cells.Item(134, 11).Value = 2992.5  # K134
$ws.Cells.Item(134, 13).Value = 2077.5  # M134


# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 784.6667  # H97
$ws.Cells.Item(97, 9).Value = 518.75  # I97
$ws.Cells.Item(97, 11).Value = 518.75  # K97
$ws.Cells.Item(97, 13).Value = -22.75  # M97


# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1093.25  # H22
$ws.Cells.Item(22, 9).Value = 1029.4  # I22
$ws.Cells.Item(22, 10).Value = 1199.6666  # J22
$ws.Cells.Item(22, 11).Value = 1029.4  # K22
$ws.Cells.Item(22, 12).Value = 1199.6666  # L22
$ws.Cells.Item(22, 13).Value = -734.4000000000001  # M22
$ws.Cells.Item(22, 14).Value = -1789.6666  # N22

$ws.Cells.Item(27, 8).Value = 1093.25  # H27
$ws.Cells.Item(27, 9).Value = 1029.4  # I27
$ws.Cells.Item(27, 10).Value = 1199.6666  # J27
$ws.Cells.Item(27, 11).Value = 1029.4  # K27
$ws.Cells.Item(27, 12).Value = 1199.6666  # L27
$ws.Cells.Item(27, 13).Value = -922.4000000000001  # M27
$ws.Cells.Item(27, 14).Value = -1413.6666  # N27

$ws.Cells.Item(42, 8).Value = 35085  # H42
$ws.Cells.Item(42, 9).Value = 0  # I42
$ws.Cells.Item(42, 10).Value = 35085  # J42
$ws.Cells.Item(42, 11).Value = 0  # K42
$ws.Cells.Item(42, 12).Value = 35085  # L42
$ws.Cells.Item(42, 13).ClearContents()  # M42
$ws.Cells.Item(42, 14).Value = -36211  # N42

$ws.Cells.Item(46, 8).Value = 750  # H46
$ws.Cells.Item(46, 9).Value = 625  # I46
$ws.Cells.Item(46, 10).Value = 1000  # J46
$ws.Cells.Item(46, 11).Value = 625  # K46
$ws.Cells.Item(46, 12).Value = 1000  # L46
$ws.Cells.Item(46, 13).Value = -437  # M46
$ws.Cells.Item(46, 14).Value = -1376  # N46

$ws.Cells.Item(49, 8).Value = 35085  # H49
$ws.Cells.Item(49, 9).Value = 0  # I49
$ws.Cells.Item(49, 10).Value = 35085  # J49
$ws.Cells.Item(49, 11).Value = 0  # K49
$ws.Cells.Item(49, 12).Value = 35085  # L49
$ws.Cells.Item(49, 13).ClearContents()  # M49
$ws.Cells.Item(49, 14).Value = -35379  # N49

$ws.Cells.Item(50, 8).Value = 29994  # H50
$ws.Cells.Item(50, 9).Value = 0  # I50
$ws.Cells.Item(50, 10).Value = 29994  # J50
$ws.Cells.Item(50, 11).Value = 0  # K50
$ws.Cells.Item(50, 12).Value = 29994  # L50
$ws.Cells.Item(50, 13).ClearContents()  # M50
$ws.Cells.Item(50, 14).Value = -31268  # N50

$ws.Cells.Item(55, 8).Value = 975.4  # H55
$ws.Cells.Item(55, 9).Value = 381.77777  # I55
$ws.Cells.Item(55, 10).Value = 1865.8334  # J55
$ws.Cells.Item(55, 11).Value = 381.77777  # K55
$ws.Cells.Item(55, 12).Value = 1865.8334  # L55
$ws.Cells.Item(55, 13).Value = -208.77777  # M55
$ws.Cells.Item(55, 14).Value = -2211.8334  # N55

$ws.Cells.Item(68, 8).Value = 2424.0833  # H68
$ws.Cells.Item(68, 9).Value = 2539.2  # I68
$ws.Cells.Item(68, 10).Value = 1848.5  # J68
$ws.Cells.Item(68, 11).Value = 2539.2  # K68
$ws.Cells.Item(68, 12).Value = 1848.5  # L68
$ws.Cells.Item(68, 13).Value = -1790.2  # M68
$ws.Cells.Item(68, 14).Value = -3346.5  # N68

$ws.Cells.Item(71, 8).Value = 2424.0833  # H71
$ws.Cells.Item(71, 9).Value = 2539.2  # I71
$ws.Cells.Item(71, 10).Value = 1848.5  # J71
$ws.Cells.Item(71, 11).Value = 12696  # K71
$ws.Cells.Item(71, 12).Value = 9242.5  # L71
$ws.Cells.Item(71, 13).Value = -8952  # M71
$ws.Cells.Item(71, 14).Value = -16730.5  # N71


# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5443.6665  # H62
$ws.Cells.Item(62, 9).Value = 5900  # I62
$ws.Cells.Item(62, 10).Value = 4873.25  # J62
$ws.Cells.Item(62, 11).Value = 5900  # K62
$ws.Cells.Item(62, 12).Value = 4873.25  # L62
$ws.Cells.Item(62, 13).Value = -5276  # M62
$ws.Cells.Item(62, 14).Value = -6121.25  # N62

$ws.Cells.Item(65, 8).Value = 5443.6665  # H65
$ws.Cells.Item(65, 9).Value = 5900  # I65
$ws.Cells.Item(65, 10).Value = 4873.25  # J65
$ws.Cells.Item(65, 11).Value = 29500  # K65
$ws.Cells.Item(65, 12).Value = 24366.25  # L65
$ws.Cells.Item(65, 13).Value = -26380  # M65
$ws.Cells.Item(65, 14).Value = -30606.25  # N65

$ws.Cells.Item(122, 8).Value = 2586.3845  # H122
$ws.Cells.Item(122, 10).Value = 2623.6667  # J122
$ws.Cells.Item(122, 12).Value = 7871.000100000001  # L122
$ws.Cells.Item(122, 14).Value = -12771.0001  # N122

$ws.Cells.Item(126, 8).Value = 2439.75  # H126
$ws.Cells.Item(126, 9).Value = 2034.625  # I126
$ws.Cells.Item(126, 11).Value = 6103.875  # K126
$ws.Cells.Item(126, 13).Value = -3633.875  # M126

